$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add a new rule row for ORGANIZATION, mirroring the existing CASE_FILE / COMPLAINT rows
$ws.Range("B22").Value = "Organization - Check participants list for NoAccess & Owner"
$ws.Range("C22").Value = "ORGANIZATION"
$ws.Range("D22").Value = "participants != null && participants.containsKey('No Access') && participants.containsKey('owner')"
$ws.Range("E22").Value = "participants['No Access'].contains(participants['owner'][0])"
$ws.Range("F22").Value = "Owners cannot be on the no-access list."

# Match formatting of the row above (row 21) - copy formats only, values already set above
$ws.Range("B21:F21").Copy()
$ws.Range("B22:F22").PasteSpecial(-4122)
$ws.Rows.Item(22).RowHeight = 45

$ws.Range("F22").Select()
